# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (G) was recalculated; write the new computed values for
# rows 2-29 (row 27 is unchanged, old value already equalled new value).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 5
    3  = 1
    4  = 3
    5  = 7
    6  = 8
    7  = 2
    8  = 4
    9  = 6
    10 = 2
    11 = 6
    12 = 4
    13 = 4
    14 = 4
    15 = 3
    16 = 5
    17 = 1
    18 = 6
    19 = 3
    20 = 2
    21 = 3
    22 = 4
    23 = 3
    24 = 1
    25 = 1
    26 = 1
    27 = 1
    28 = 4
    29 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
